$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells (AD1:AF1) with the same style as the existing headers ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting (bold/border/alignment) from an existing header cell
# so the new header cells pick up the same style index.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Populate the season record (Wins/Losses/Ties) for every data row ---
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 67
    $ws.Cells.Item($r, 31).Value = 95
    $ws.Cells.Item($r, 32).Value = 0
}
